$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a value as text (keeping default/no-style look) for the
# year_max columns, which store years as text (e.g. "2024") rather than
# numbers. Applying "@" before the write keeps Excel from auto-converting
# the numeric-looking string to a number; resetting the style back to
# "Normal" afterwards avoids leaving a stray number-format style behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("N2") "2024"
$ws.Range("P2").Value = 4.1
Set-TextValue $ws.Range("Q2") "2024"
$ws.Range("S2").Value = 12.60752813494749
Set-TextValue $ws.Range("T2") "2024"
$ws.Range("V2").Value = 18.00322580645161
Set-TextValue $ws.Range("W2") "2024"
$ws.Range("Y2").Value = 8.480645161290322
Set-TextValue $ws.Range("AC2") "2024"
$ws.Range("AE2").Value = 8.525690214803117
